# "2e kans - 3 Lotte Houwen"
# Rewrites the test-case rows of the test script with a new (2nd-chance)
# round of test steps, tweaks the date/time header cells, widens column G,
# and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header area: date of test run changes, and the "time" cell becomes
#     free text instead of a real time value -----------------------------
$ws.Range("G7").Value = "2/28/2018"
$ws.Range("G8").Value = "20:0000"

# --- Test-case rows ------------------------------------------------------

# Row 10 - bestand inlezen: virushostdb.tsv meegeven
$ws.Range("B10").Value = "bestand inlezen: virushostdb.tsv meegeven"
$ws.Range("C10").Value = "Geen exception. Bestand is in te lezen en kan nu gefilterd en gesorteerd gaan worden."
$ws.Range("D10").Value = "x"
$ws.Rows.Item(10).RowHeight = 25.5

# Row 11 - Bestand dmv typen invoeren in bestandveld
$ws.Range("B11").Value = "Bestand dmv typen invoeren in bestandveld"
$ws.Range("C11").Value = "Er gebeurd niets"
$ws.Range("D11").Value = "x"

# Row 12 - Filteren/Sorteren ssRNA virusses, zelfde hosts
$ws.Range("B12").Value = "Filteren/Sorteren: bij classe ssRNA virusses en bij host 9796 en 9796 (dezelfde hosts)"
$ws.Range("C12").Value = "dezelfde lijsten terug, op dezelfde volgorde en bij overeenkomst alle virussen terug die bij de 2 tekstarea's staan."
$ws.Range("D12").Value = "x"
$ws.Range("G12").Value = "sorteeroptie virusid"

# Row 13 - Filteren/Sorteren ssRNA virusses, andere host
$ws.Range("B13").Value = "Filteren/Sorteren: bij classe ssRNA virusses en bij host 9796 en 9790 (dezelfde hosts)"
$ws.Range("C13").Value = "In elk tekstarea een eigen lijst van virusid bij overeenkomst 1 overeenkomstige virusid"
$ws.Range("D13").Value = "x"
$ws.Range("G13").Value = "sorteeroptie virusid"
$ws.Rows.Item(13).RowHeight = 30.75

# Row 14 - fully cleared (was test step 5, now blank pending new content)
# NOTE: ClearContents() is unreliable on cells that anchor a merged range in
# this runtime, so every clear below uses Value = "" instead, which is
# honoured consistently (and yields the same clean empty <c> element).
$ws.Range("A14").Value = ""
$ws.Range("B14").Value = ""
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""

# Rows 15-20 - text content cleared; the data/format columns (E/F/G/H) were
# already empty. Row A-numbers for these rows stay as-is.
$ws.Range("B15").Value = ""
$ws.Range("C15").Value = ""
$ws.Range("D15").Value = ""

$ws.Range("B16").Value = ""
$ws.Range("C16").Value = ""
$ws.Range("D16").Value = ""
$ws.Rows.Item(16).RowHeight = 12.75

$ws.Range("B17").Value = ""
$ws.Range("C17").Value = ""
$ws.Range("D17").Value = ""
$ws.Rows.Item(17).RowHeight = 12.75

$ws.Range("B18").Value = ""
$ws.Range("C18").Value = ""
$ws.Range("D18").Value = ""
$ws.Rows.Item(18).RowHeight = 12.75

$ws.Range("B19").Value = ""
$ws.Range("C19").Value = ""
$ws.Range("E19").Value = ""
$ws.Range("G19").Value = ""
$ws.Rows.Item(19).RowHeight = 12.75

$ws.Range("B20").Value = ""
$ws.Range("C20").Value = ""
$ws.Range("D20").Value = ""
$ws.Rows.Item(20).RowHeight = 12.75

# --- Column G a touch wider so the new text fits -------------------------
$ws.Columns.Item(7).ColumnWidth = 9.6

# --- Move the active selection -------------------------------------------
$ws.Range("L14").Select()
